$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "...Képek és videok keresése a teko weboldaláról..." paragraph.
#   - "videok" -> "videók"
#   - "teko"   -> "teko.hu"
#   - drop the (now stale) spell-check proofErr markers around those words
#   - the "_GoBack" bookmark moves to just after "teko.hu"
# ---------------------------------------------------------------------------

$p1 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -match "11\.18-11:30-14:17") {
        $p1 = $d.Paragraphs($i)
        break
    }
}

$findRange = $p1.Range.Duplicate
$findRange.Find.ClearFormatting()
$findRange.Find.Execute(
    "Képek és videok keresése a teko weboldaláról",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Képek és videók keresése a teko.hu weboldaláról", 2) | Out-Null

# Locate the freshly written "teko.hu" so the bookmark can be re-anchored
# right after it.
$tekoRange = $p1.Range.Duplicate
$tekoRange.Find.ClearFormatting()
$tekoRange.Find.Execute("teko.hu") | Out-Null
$bookmarkPos = $tekoRange.End

$oldBookmark = $d.Bookmarks("_GoBack")
$oldBookmark.Delete()
$d.Bookmarks.Add("_GoBack", $d.Range($bookmarkPos, $bookmarkPos)) | Out-Null

# ---------------------------------------------------------------------------
# Change 2: "11.18-20:50-" paragraph gains the rest of that log entry plus a
# new trailing sentence about resizing behaviour.
# ---------------------------------------------------------------------------

$p2 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -eq "11.18-20:50-`r") {
        $p2 = $d.Paragraphs($i)
        break
    }
}

$tailRange = $p2.Range.Duplicate
$tailRange.Find.ClearFormatting()
$tailRange.Find.Execute("-20:50-") | Out-Null
$tailRange.Text = "-20:50-20:20 – A videó control panelje felül volt a menübe amikor ki volt nyitva (Kijavítva)."

$insertPos = $p2.Range.End - 1
$d.Range($insertPos, $insertPos).InsertAfter(" Átméretezések mind gépi és mobilos felületen") | Out-Null
